# Update main GSC export data (Breadcrumbs "Chart" sheet):
#  - Drop the oldest day (2025-10-21), shifting every later day/row up by one.
#  - Append a new trailing day (2026-01-19) with its Invalid/Valid counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Remove the row for 2025-10-21 (row 2); this shifts all subsequent rows
# up by one and drops the now-unused shared string automatically.
$ws.Rows.Item(2).Delete()

# The previously-last row (91) is now empty after the shift; populate it
# with the new trailing date and its counts. Set the value with a leading
# apostrophe so Excel stores it as literal text (matching the other date
# cells) instead of auto-converting the yyyy-MM-dd-looking text to a date
# serial number. Re-applying the format from an existing date cell clears
# the "quote prefix" flag so the new cell keeps the same default style as
# its neighbors.
$ws.Range("A91").Value = "'2026-01-19"
$ws.Range("A2").Copy()
$ws.Range("A91").PasteSpecial(-4122)

$ws.Range("B91").Value = 0
$ws.Range("C91").Value = 26
